# Insurance Understanding.xlsx - apply commit "New file and modification to Ins Understanding"
#
# Summary of the edit:
#  - Sheet1: row1 data shifted one column left (B:D -> A:C); the small
#    "regi / CS INS / process understanding" block shifted from M7:O9
#    into L6:M9; view scrolled/selected to D14.
#  - TLI: selection moved to F2.
#  - CI: selection moved to B46 and becomes the active sheet/tab.
#  - HI: no longer the active tab.
#  - Two new sheets appended at the end: "Sheet2" and "Sheet3", each with a
#    small table of data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet1 content rework
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 1: shift the three labels one column to the left (B1:D1 -> A1:C1).
$ws1.Range("A1").Value = "30 days"
$ws1.Range("B1").Value = "10 steps"
$ws1.Range("C1").Value = "roughly 3 days for each"
$ws1.Range("D1").ClearContents()

# Rows 6-9: "regi" moves from M7 up to L6; "CS INS" moves from N7 to L7;
# "process understanding " moves from O7 to M7; the O8/O9 comments move to
# M8/M9 respectively.
$ws1.Range("L6").Value = "regi"
$ws1.Range("L7").Value = "CS INS"
$ws1.Range("M7").Value = "process understanding "
$ws1.Range("N7").ClearContents()
$ws1.Range("O7").ClearContents()
$ws1.Range("M8").Value = "link validity"
$ws1.Range("O8").ClearContents()
$ws1.Range("M9").Value = "cancel reason codes & when assigned"
$ws1.Range("O9").ClearContents()

# ---------------------------------------------------------------------
# 2. New sheets "Sheet2" and "Sheet3", appended after HI.
# ---------------------------------------------------------------------
$hi = $wb.Worksheets.Item("HI")

$sheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $hi)
$sheet2.Name = "Sheet2"

$sheet2.Range("A1").Value = "Dimensions at eligibiity Page"
$sheet2.Range("B1").Value = "CI"
$sheet2.Range("C1").Value = "HI"
$sheet2.Range("D1").Value = "TLI"
$sheet2.Range("E1").Value = "Dimension/Facr"

$sheet2.Range("A2").Value = "Age"
$sheet2.Range("B2").Value = "YES"
$sheet2.Range("F2").Value = "Date of "

$sheet2.Range("A3").Value = "Smoke/"
$sheet2.Range("F3").Value = "na"
$sheet2.Range("G3").Value = "Raise Reqest"

$sheet2.Range("A4").Value = "Applicant City"
$sheet2.Range("I4").Value = "uniquw"

$sheet2.Range("A5").Value = "Gross Mn"

$sheet2.Range("A6").Value = "x"
$sheet2.Range("D6").Value = "YES"

$sheet2.Range("A14").Value = "Appl Dimension"

$sheet2.Range("A15").Value = "First"

$sheet2.Columns.Item(1).ColumnWidth = 27.140625

$sheet3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet2)
$sheet3.Name = "Sheet3"

$sheet3.Range("A1").Value = "No Visitrs"
$sheet3.Range("B1").Value = "WH"
$sheet3.Range("C1").Value = "Visits"
$sheet3.Range("D1").Value = "FieldSession"
$sheet3.Range("E1").Value = "Sum"
$sheet3.Range("F1").Value = "Filter Condition - N"
$sheet3.Range("I1").Value = "Finance"
$sheet3.Range("J1").Value = "TI"
$sheet3.Range("K1").Value = "Ops"

$sheet3.Range("A2").Value = "No Offers"
$sheet3.Range("I2").Value = "Finance"
$sheet3.Range("J2").Value = "TI"
$sheet3.Range("K2").Value = "Ops"

$sheet3.Range("A3").Value = "No Searches"
$sheet3.Range("I3").Value = "Finance"
$sheet3.Range("J3").Value = "TI"
$sheet3.Range("K3").Value = "Ops"

$sheet3.Range("A4").Value = "N "
$sheet3.Range("J4").Value = "TI"

$sheet3.Range("K5").Value = "Ops"

$sheet3.Columns.Item(1).ColumnWidth = 11.85546875
$sheet3.Columns.Item(4).ColumnWidth = 14.5703125
$sheet3.Columns.Item(6).ColumnWidth = 15.85546875

# ---------------------------------------------------------------------
# 3. View / selection state per sheet.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("D14").Select()

$tli = $wb.Worksheets.Item("TLI")
$tli.Activate()
$tli.Range("F2").Select()

$sheet2.Activate()
$sheet2.Range("A6").Select()

$sheet3.Activate()
$sheet3.Range("D1").Select()

# CI becomes the final active / selected sheet (matches activeTab=2 and
# the tabSelected="1" that moves off HI and onto CI).
$ci = $wb.Worksheets.Item("CI")
$ci.Activate()
$ci.Range("B46").Select()
